$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7: PN Records ---------------------------------------------------
# Set URL before Name so the shared-string table picks up the URL first
# (matches the order new strings were authored in the target workbook).
$ws.Cells.Item(7, 3).Value = "http://192.168.100.19/thaimaiapp/api/mother/mPN_Record"
$ws.Cells.Item(7, 2).Value = "PN Records"
$ws.Cells.Item(7, 4).Value = "picmeId=1000000000001" + [char]10 + "mid=1"

# Style B7/C7 like the other "name"/"url" cells in the table (center/center).
$ws.Cells.Item(7, 2).HorizontalAlignment = -4108
$ws.Cells.Item(7, 2).VerticalAlignment = -4108
$ws.Cells.Item(7, 3).HorizontalAlignment = -4108
$ws.Cells.Item(7, 3).VerticalAlignment = -4108

# D7 uses the wrap-text centered style (same as D6).
$ws.Cells.Item(7, 4).HorizontalAlignment = -4108
$ws.Cells.Item(7, 4).VerticalAlignment = -4108
$ws.Cells.Item(7, 4).WrapText = $true

$ws.Rows.Item(7).RowHeight = 30

# --- Row 8: Delivery Insert ----------------------------------------------
$ws.Cells.Item(8, 1).Value = "POST"
$ws.Cells.Item(8, 3).Value = "http://192.168.100.19/thaimaiapp/api/mother/mDeleveryDetailsInsert"
$ws.Cells.Item(8, 2).Value = "Delivery Insert"
$ws.Cells.Item(8, 4).Value = "picmeId=1000000000001" + [char]10 + "mid=1"

$ws.Cells.Item(8, 2).HorizontalAlignment = -4108
$ws.Cells.Item(8, 2).VerticalAlignment = -4108

# C8 gets its own (slightly different) centered style/font.
$ws.Cells.Item(8, 3).Font.Name = "Calibri"
$ws.Cells.Item(8, 3).HorizontalAlignment = -4108
$ws.Cells.Item(8, 3).VerticalAlignment = -4108

$ws.Cells.Item(8, 4).HorizontalAlignment = -4108
$ws.Cells.Item(8, 4).VerticalAlignment = -4108
$ws.Cells.Item(8, 4).WrapText = $true

$ws.Rows.Item(8).RowHeight = 30

# --- Column C width (was best-fit, now a fixed custom width) -------------
$ws.Columns.Item(3).ColumnWidth = 66

# --- Selection matches the post-edit cursor position ----------------------
$ws.Range("E8").Select()
